$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'62.620.41"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = "'  -2.88%  "
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.Value = "'3.187.93"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = "'  -3.63%  "
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.Value = "'0.998"
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.Value = "'  -0.27%  "
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.Value = "'508.98"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = "'  -4.23%  "
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.Value = "'167.45"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = "'  -7.95%  "
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.Value = "'0.573"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = "'  -5.43%  "
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.Value = "'0.998"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = "'  -0.23%  "
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.Value = "'3.183.32"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.Value = "'  -3.68%  "
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.Value = "'0.582"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = "'  -6.31%  "
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.Value = "'51.07"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = "'  -13.44%  "
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.Value = "'0.128"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.Value = "'  -4.37%  "
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.Value = "'0.0000248"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = "'  -6.07%  "
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.Value = "'8.65"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = "'  -5.18%  "
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.Value = "'3.732.60"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = "'  -2.45%  "
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.Value = "'3.207.20"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = "'  -2.87%  "
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.Value = "'0.114"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = "'  -2.83%  "
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.Value = "'62.244.28"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = "'  -3.14%  "
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.Value = "'16.79"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = "'  -5.56%  "
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.Value = "'10.77"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = "'  -3.27%  "
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.Value = "'0.923"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = "'  -4.04%  "
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.Value = "'364.02"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = "'  -2.91%  "
$cell.Style = "Normal"

$cell = $ws.Range("E23")
$cell.Value = "'  +4.56%  "
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.Value = "'10.97"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = "'  -2.28%  "
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.Value = "'79.03"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = "'  -2.75%  "
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.Value = "'3.53"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = "'  -7.84%  "
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.Value = "'6.08"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = "'  -0.08%  "
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.Value = "'2.59"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.Value = "'  -3.50%  "
$cell.Style = "Normal"

$cell = $ws.Range("D29")
$cell.Value = "'10.82"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.Value = "'  -6.92%  "
$cell.Style = "Normal"

$cell = $ws.Range("D30")
$cell.Value = "'7.83"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.Value = "'  -7.30%  "
$cell.Style = "Normal"

$cell = $ws.Range("D31")
$cell.Value = "'27.67"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.Value = "'  -4.74%  "
$cell.Style = "Normal"

$cell = $ws.Range("D32")
$cell.Value = "'610.66"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.Value = "'  -5.10%  "
$cell.Style = "Normal"

$cell = $ws.Range("D33")
$cell.Value = "'6.37"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.Value = "'  -6.36%  "
$cell.Style = "Normal"

$cell = $ws.Range("D34")
$cell.Value = "'10.84"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.Value = "'  -4.50%  "
$cell.Style = "Normal"

$cell = $ws.Range("D35")
$cell.Value = "'0.101"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.Value = "'  -3.78%  "
$cell.Style = "Normal"

$cell = $ws.Range("D36")
$cell.Value = "'56.26"
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.Value = "'  -5.99%  "
$cell.Style = "Normal"

$cell = $ws.Range("D37")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.Value = "'  +0.08%  "
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.Value = "'34.43"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.Value = "'  -6.78%  "
$cell.Style = "Normal"

$cell = $ws.Range("D39")
$cell.Value = "'0.363"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.Value = "'  -8.22%  "
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.Value = "'0.999"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = "'  +0.09%  "
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.Value = "'0.0₃0693"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = "'  -5.72%  "
$cell.Style = "Normal"

$cell = $ws.Range("B42")
$cell.Value = "'Fetch.AI"
$cell.Style = "Normal"
$cell = $ws.Range("C42")
$cell.Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.Value = "'2.52"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = "'  +1.04%  "
$cell.Style = "Normal"

$cell = $ws.Range("B43")
$cell.Value = "'Kaspa"
$cell.Style = "Normal"
$cell = $ws.Range("C43")
$cell.Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.Value = "'0.120"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = "'  -5.99%  "
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.Value = "'2.788.10"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = "'  -4.08%  "
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.Value = "'2.92"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = "'  -0.30%  "
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.Value = "'2.64"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = "'  -0.93%  "
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.Value = "'0.0380"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = "'  -6.00%  "
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.Value = "'2.49"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = "'  -7.93%  "
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.Value = "'2.92"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = "'  -1.20%  "
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.Value = "'134.28"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = "'  +0.91%  "
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.Value = "'0.121"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = "'  -4.79%  "
$cell.Style = "Normal"
